$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, new nombre_aides (C), new montant_total (D)
$updates = @(
    @(2, 322380, 410843293),
    @(3, 259, 309979),
    @(4, 317, 453707),
    @(8, 866, 1273795),
    @(10, 117398, 172022519),
    @(11, 143, 211513),
    @(12, 59711, 86184648),
    @(16, 4009, 5688992),
    @(19, 71, 103035),
    @(20, 6725, 9388533),
    @(22, 77849, 97046177),
    @(27, 288, 413147),
    @(28, 32564, 47669242),
    @(30, 11529, 16584949),
    @(32, 4, 6000),
    @(33, 1560, 2191307),
    @(35, 1844, 2603270),
    @(36, 97589, 122794547),
    @(37, 67, 77868),
    @(42, 903, 1329185),
    @(44, 44481, 65188847),
    @(46, 9184, 13178380),
    @(48, 1408, 1954603),
    @(51, 2325, 3245929),
    @(52, 69352, 86989444),
    @(53, 42, 47883),
    @(57, 383, 562730),
    @(59, 28284, 41479994),
    @(62, 11178, 16164170),
    @(64, 1360, 1902289),
    @(68, 1487, 2082584),
    @(70, 20592, 26975819),
    @(74, 7616, 11153008),
    @(76, 5148, 7474910),
    @(77, 494, 699739),
    @(78, 280, 393083),
    @(79, 141494, 176410978),
    @(80, 70, 84321),
    @(85, 63711, 93377216),
    @(88, 29867, 43207116),
    @(91, 2855, 4037318),
    @(92, 33426, 45305052),
    @(95, 28, 40314),
    @(96, 8086, 11888715),
    @(98, 7442, 10798839),
    @(100, 538, 763651),
    @(101, 499, 719891),
    @(102, 10098, 15183364),
    @(104, 2515, 4043791),
    @(106, 3380, 5412425),
    @(108, 154, 248320),
    @(109, 192, 288468),
    @(110, 142213, 175859022),
    @(114, 951, 1395536),
    @(116, 52884, 77518194),
    @(118, 27230, 39452525),
    @(119, 1313, 1795784),
    @(122, 2282, 3206569),
    @(124, 516974, 682834769),
    @(126, 213, 313509),
    @(129, 1381, 2046986),
    @(131, 209127, 307433450),
    @(132, 405, 604250),
    @(134, 184300, 268004291),
    @(135, 4, 6000),
    @(137, 2854, 4008969),
    @(139, 6420, 9069280),
    @(142, 44941, 60005115),
    @(143, 29, 40262),
    @(146, 5, 7500),
    @(148, 14137, 20728558),
    @(149, 3787, 5461520),
    @(154, 387, 547163),
    @(155, 17730, 23436553),
    @(159, 7222, 10505701),
    @(161, 5047, 7264413),
    @(163, 281, 388931),
    @(164, 270, 386364),
    @(166, 18005, 28860627),
    @(167, 1985, 3214449),
    @(168, 266, 425661),
    @(170, 63, 107690),
    @(171, 97, 160449),
    @(172, 88166, 110225791),
    @(179, 33975, 49824677),
    @(181, 13097, 18924462),
    @(185, 1669, 2346650),
    @(187, 239691, 297920989),
    @(188, 136, 146087),
    @(193, 882, 1297345),
    @(195, 86861, 127329295),
    @(198, 33198, 47789043),
    @(201, 5121, 7290277),
    @(204, 4926, 6820989),
    @(207, 265303, 328314392),
    @(209, 255, 365087),
    @(211, 11, 12731),
    @(216, 95461, 139656018),
    @(217, 90, 134199),
    @(219, 51710, 74743148),
    @(222, 4684, 6577031),
    @(225, 5824, 8057434),
    @(228, 106890, 133663492),
    @(235, 49670, 72762868),
    @(237, 12494, 17964997),
    @(241, 2539, 3552972),
    @(242, 259193, 327300068),
    @(244, 251, 360312),
    @(248, 833, 1223563),
    @(250, 96169, 140916287),
    @(253, 65459, 94866547),
    @(255, 2417, 3408238),
    @(258, 4656, 6541560)
)

foreach ($u in $updates) {
    $row = $u[0]
    $ws.Cells.Item($row, 3).Value = $u[1]
    $ws.Cells.Item($row, 4).Value = $u[2]
}

$wb.Save()